$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'puma tights for women'
$ws.Range("A2").Value = 'puma tights women'
$ws.Range("A3").Value = 'puma x'
$ws.Range("A4").Value = 'purple compression shorts'
$ws.Range("A5").Value = 'purple plus size stockings'
$ws.Range("A6").Value = 'purple shorts women high waist'
$ws.Range("A7").Value = 'purple tights plus size'
$ws.Range("A8").Value = 'purple under armour shorts'
$ws.Range("A9").Value = 'putting plane'
$ws.Range("A10").Value = 'quad bike games'
$ws.Range("A11").Value = 'quad compression'
$ws.Range("A12").Value = 'quad pants'
$ws.Range("A13").Value = 'quad stretch'
$ws.Range("A14").Value = 'quad stretch strap'
$ws.Range("A15").Value = 'quadricep compression'
$ws.Range("A16").Value = 'quick dry capris women'
$ws.Range("A17").Value = 'quick dry hiking capris women'
$ws.Range("A18").Value = 'race clothes for women'
$ws.Range("A19").Value = 'rainbeau curves compression'
$ws.Range("A20").Value = 'rainbow apparel'
$ws.Range("A21").Value = 'rainbow apperal'
$ws.Range("A22").Value = 'rainbow athletic tape'
$ws.Range("A23").Value = 'rainbow bike'
$ws.Range("A24").Value = 'rainbow black shorts'
$ws.Range("A25").Value = 'rainbow bright'
$ws.Range("A26").Value = 'rainbow bright clothes'
$ws.Range("A27").Value = 'rainbow capri'
$ws.Range("A28").Value = 'rainbow clothes women'
$ws.Range("A29").Value = 'rainbow clothing for women'
$ws.Range("A30").Value = 'rainbow clothing men'
$ws.Range("A31").Value = 'rainbow clothing women'
$ws.Range("A32").Value = 'rainbow elastic string'
$ws.Range("A33").Value = 'rainbow gear'
$ws.Range("A34").Value = 'rainbow gym shorts women'
$ws.Range("A35").Value = 'rainbow hot pants'
$ws.Range("A36").Value = 'rainbow kinesiology tape'
$ws.Range("A37").Value = 'rainbow leg warmers for women'
$ws.Range("A38").Value = 'rainbow legging'
$ws.Range("A39").Value = 'rainbow legginga'
$ws.Range("A40").Value = 'rainbow leggings'
$ws.Range("A41").Value = 'rainbow leggings for women'
$ws.Range("A42").Value = 'rainbow leggings women'
$ws.Range("A43").Value = 'rainbow long live'
$ws.Range("A44").Value = 'rainbow lycra'
$ws.Range("A45").Value = 'rainbow pants'
$ws.Range("A46").Value = 'rainbow pants men'
$ws.Range("A47").Value = 'rainbow pastel'
$ws.Range("A48").Value = 'rainbow pocket chart'
$ws.Range("A49").Value = 'rainbow reflective leggings'
$ws.Range("A50").Value = 'rainbow reflective tape'
$ws.Range("A51").Value = 'rainbow road'
$ws.Range("A52").Value = 'rainbow run'
$ws.Range("A53").Value = 'rainbow runner'
$ws.Range("A54").Value = 'rainbow runners'
$ws.Range("A55").Value = 'rainbow running'
$ws.Range("A56").Value = 'rainbow running shorts'
$ws.Range("A57").Value = 'rainbow sets women clothing'
$ws.Range("A58").Value = 'rainbow shop'
$ws.Range("A59").Value = 'rainbow short'
$ws.Range("A60").Value = 'rainbow short shorts'
$ws.Range("A61").Value = 'rainbow shorts'
$ws.Range("A62").Value = 'rainbow shorts for women'
$ws.Range("A63").Value = 'rainbow shorts women'
$ws.Range("A64").Value = 'rainbow spandex'
$ws.Range("A65").Value = 'rainbow spandex fabric'
$ws.Range("A66").Value = 'rainbow stockings plus size'
$ws.Range("A67").Value = 'rainbow store clothes for women'
$ws.Range("A68").Value = 'rainbow stripe leggings'
$ws.Range("A69").Value = 'rainbow tight'
$ws.Range("A70").Value = 'rainbow tights'
$ws.Range("A71").Value = 'rainbow tights for women'
$ws.Range("A72").Value = 'rainbow underarmour'
$ws.Range("A73").Value = 'rainbow web 3'
$ws.Range("A74").Value = 'rainbow week'
$ws.Range("A75").Value = 'rainbow women'
$ws.Range("A76").Value = 'rainbow women shorts'
$ws.Range("A77").Value = 'rainbow women top'
$ws.Range("A78").Value = 'rainbow womens'
$ws.Range("A79").Value = 'rainbow womens clothes'
$ws.Range("A80").Value = 'rainbow womens clothing'
$ws.Range("A81").Value = 'rainbow womens leg warmers'
$ws.Range("A82").Value = 'rainbow womens tights'
$ws.Range("A83").Value = 'rainbow workout clothes'
$ws.Range("A84").Value = 'rainbow workout leggings'
$ws.Range("A85").Value = 'rainbow workout pants'
$ws.Range("A86").Value = 'rainbow yoga'
$ws.Range("A87").Value = 'rainbow yoga pants'
$ws.Range("A88").Value = 'rainbows womens'
$ws.Range("A89").Value = 'range performance'
$ws.Range("A90").Value = 'range runners'
$ws.Range("A91").Value = 'reading goal'
$ws.Range("A92").Value = 'real women drive trucks'
$ws.Range("A93").Value = 'rebook basketball'
$ws.Range("A94").Value = 'record runner'
$ws.Range("A95").Value = 'recover post workout'
$ws.Range("A96").Value = 'recovery compression'
$ws.Range("A97").Value = 'recovery compression leggings'
$ws.Range("A98").Value = 'recovery compression leggings women'
$ws.Range("A99").Value = 'recovery compression pants'
$ws.Range("A100").Value = 'recovery compression pants men'
